$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues; bakes formula result into a literal value, preserving cell style
}

# Promote rows that become highlighted ("whale mover") rows by copying the
# format from row 15, which already uses the highlight style (s=4/5/6/7).
$ws.Range("A15:G15").Copy() | Out-Null
$ws.Range("A6:G6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A18:G18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A24:G24").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A26:G26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A44:G44").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 2
$ws.Cells.Item(2,1).Value = 1
Set-TextValue $ws.Cells.Item(2,2) '0xc3d8f7d9a88ed1a7dcc3414289a42956a8fd4941'
Set-TextValue $ws.Cells.Item(2,3) '761,082,598,387,042'
Set-TextValue $ws.Cells.Item(2,4) '761 Trillion'
Set-TextValue $ws.Cells.Item(2,5) '0.7611% '

# Row 3
$ws.Cells.Item(3,1).Value = 2
Set-TextValue $ws.Cells.Item(3,2) '0x28e87335c43b9b28d658b6825811744c251f5974'
Set-TextValue $ws.Cells.Item(3,3) '756,107,607,679,423'
Set-TextValue $ws.Cells.Item(3,4) '756 Trillion'
Set-TextValue $ws.Cells.Item(3,5) '0.7561% '

# Row 4
$ws.Cells.Item(4,1).Value = 3
Set-TextValue $ws.Cells.Item(4,2) '0x110d4a2fb45e361194a17df10635efdb41d7579c'
Set-TextValue $ws.Cells.Item(4,3) '603,980,172,214,805'
Set-TextValue $ws.Cells.Item(4,4) '604 Trillion'
Set-TextValue $ws.Cells.Item(4,5) '0.6040% '

# Row 5
$ws.Cells.Item(5,1).Value = 4
Set-TextValue $ws.Cells.Item(5,2) '0x573068628665a40aa2dddfd5ec3c90167424a9ee'
Set-TextValue $ws.Cells.Item(5,3) '602,821,058,845,833'
Set-TextValue $ws.Cells.Item(5,4) '603 Trillion'
Set-TextValue $ws.Cells.Item(5,5) '0.6028% '

# Row 6
$ws.Cells.Item(6,1).Value = 5
Set-TextValue $ws.Cells.Item(6,2) '0xcdc162e71e7517e94c0c72de16ca6315349d47fc'
Set-TextValue $ws.Cells.Item(6,3) '529,250,426,560,027'
Set-TextValue $ws.Cells.Item(6,4) '529 Trillion'
Set-TextValue $ws.Cells.Item(6,5) '0.5293% '
Set-TextValue $ws.Cells.Item(6,6) '29 Trillion More'
Set-TextValue $ws.Cells.Item(6,7) '28,955,311,081,981'

# Row 7
$ws.Cells.Item(7,1).Value = 6
Set-TextValue $ws.Cells.Item(7,2) '0x44cce737956e0774ed45fc88efbc733c0320ac02'
Set-TextValue $ws.Cells.Item(7,3) '524,355,454,029,402'
Set-TextValue $ws.Cells.Item(7,4) '524 Trillion'
Set-TextValue $ws.Cells.Item(7,5) '0.5244% '

# Row 8
$ws.Cells.Item(8,1).Value = 7
Set-TextValue $ws.Cells.Item(8,2) '0x7fcc2410d49222563a5cd46c0ceaffb72cf195d4'
Set-TextValue $ws.Cells.Item(8,3) '432,594,179,759,658'
Set-TextValue $ws.Cells.Item(8,4) '433 Trillion'
Set-TextValue $ws.Cells.Item(8,5) '0.4326% '

# Row 9
$ws.Cells.Item(9,1).Value = 8
Set-TextValue $ws.Cells.Item(9,2) '0xac3e1d277c3442ac1c695a369d87004aed308a7d'
Set-TextValue $ws.Cells.Item(9,3) '419,994,932,243,674'
Set-TextValue $ws.Cells.Item(9,4) '420 Trillion'
Set-TextValue $ws.Cells.Item(9,5) '0.4200% '

# Row 10
$ws.Cells.Item(10,1).Value = 9
Set-TextValue $ws.Cells.Item(10,2) '0xbd09e0594fbdbc5f73fe5db01bdc3bc2a19ec2d7'
Set-TextValue $ws.Cells.Item(10,3) '419,360,533,364,969'
Set-TextValue $ws.Cells.Item(10,4) '419 Trillion'
Set-TextValue $ws.Cells.Item(10,5) '0.4194% '

# Row 11
$ws.Cells.Item(11,1).Value = 10
Set-TextValue $ws.Cells.Item(11,2) '0x0932767b51f1faedc5a86fb935eec8032f62eea8'
Set-TextValue $ws.Cells.Item(11,3) '415,228,348,905,327'
Set-TextValue $ws.Cells.Item(11,4) '415 Trillion'
Set-TextValue $ws.Cells.Item(11,5) '0.4152% '

# Row 12
$ws.Cells.Item(12,1).Value = 11
Set-TextValue $ws.Cells.Item(12,2) '0x7167d70b2e6b167633356254bf22f5d32aedcd5b'
Set-TextValue $ws.Cells.Item(12,3) '405,513,667,501,757'
Set-TextValue $ws.Cells.Item(12,4) '406 Trillion'
Set-TextValue $ws.Cells.Item(12,5) '0.4055% '

# Row 13
$ws.Cells.Item(13,1).Value = 12
Set-TextValue $ws.Cells.Item(13,2) '0x7395cb62e4405b6c1174c2329f444af6ee7bdfd9'
Set-TextValue $ws.Cells.Item(13,3) '362,521,477,411,729'
Set-TextValue $ws.Cells.Item(13,4) '363 Trillion'
Set-TextValue $ws.Cells.Item(13,5) '0.3625% '

# Row 14
$ws.Cells.Item(14,1).Value = 13
Set-TextValue $ws.Cells.Item(14,2) '0xf7625178ef07107edc2005a0e5d2fc411573f381'
Set-TextValue $ws.Cells.Item(14,3) '357,037,237,086,702'
Set-TextValue $ws.Cells.Item(14,4) '357 Trillion'
Set-TextValue $ws.Cells.Item(14,5) '0.3570% '

# Row 15
$ws.Cells.Item(15,1).Value = 14
Set-TextValue $ws.Cells.Item(15,2) '0x69fe97ce030074b37cbaf3ee46e9f68ca8712099'
Set-TextValue $ws.Cells.Item(15,3) '324,017,179,896,078'
Set-TextValue $ws.Cells.Item(15,4) '324 Trillion'
Set-TextValue $ws.Cells.Item(15,5) '0.3240% '
Set-TextValue $ws.Cells.Item(15,6) '83 Trillion Less'
Set-TextValue $ws.Cells.Item(15,7) '82,728,944,370,975'

# Row 16
$ws.Cells.Item(16,1).Value = 15
Set-TextValue $ws.Cells.Item(16,2) '0xd96622a9099d758f8d6664ae702a59e9d548ed23'
Set-TextValue $ws.Cells.Item(16,3) '263,241,099,413,936'
Set-TextValue $ws.Cells.Item(16,4) '263 Trillion'
Set-TextValue $ws.Cells.Item(16,5) '0.2632% '

# Row 17
$ws.Cells.Item(17,1).Value = 16
Set-TextValue $ws.Cells.Item(17,2) '0x924e8fc81484781b8057db784266017fce1af136'
Set-TextValue $ws.Cells.Item(17,3) '250,088,171,286,383'
Set-TextValue $ws.Cells.Item(17,4) '250 Trillion'
Set-TextValue $ws.Cells.Item(17,5) '0.2501% '

# Row 18
$ws.Cells.Item(18,1).Value = 17
Set-TextValue $ws.Cells.Item(18,2) '0x6111742b05e8ece0d9fedb82bcfdc597be7b43d9'
Set-TextValue $ws.Cells.Item(18,3) '232,191,361,549,013'
Set-TextValue $ws.Cells.Item(18,4) '232 Trillion'
Set-TextValue $ws.Cells.Item(18,5) '0.2322% '
Set-TextValue $ws.Cells.Item(18,6) '94 Trillion More'
Set-TextValue $ws.Cells.Item(18,7) '94,248,199,228,736'

# Row 19
$ws.Cells.Item(19,1).Value = 18
Set-TextValue $ws.Cells.Item(19,2) '0x14dd63527333ff7269a59985658224f96faab068'
Set-TextValue $ws.Cells.Item(19,3) '217,372,420,035,141'
Set-TextValue $ws.Cells.Item(19,4) '217 Trillion'
Set-TextValue $ws.Cells.Item(19,5) '0.2174% '

# Row 20
$ws.Cells.Item(20,1).Value = 19
Set-TextValue $ws.Cells.Item(20,2) '0xfdd50de023c9a705d9086bf821d15c7450ee93bf'
Set-TextValue $ws.Cells.Item(20,3) '206,708,435,907,937'
Set-TextValue $ws.Cells.Item(20,4) '207 Trillion'
Set-TextValue $ws.Cells.Item(20,5) '0.2067% '

# Row 21
$ws.Cells.Item(21,1).Value = 20
Set-TextValue $ws.Cells.Item(21,2) '0x7146f34d166379b4ab5220f5eefd7c79835a3c04'
Set-TextValue $ws.Cells.Item(21,3) '195,405,153,664,131'
Set-TextValue $ws.Cells.Item(21,4) '195 Trillion'
Set-TextValue $ws.Cells.Item(21,5) '0.1954% '

# Row 22
$ws.Cells.Item(22,1).Value = 21
Set-TextValue $ws.Cells.Item(22,2) '0x7b11f31fc0d0a79717ec025d411ac5e899ac7116'
Set-TextValue $ws.Cells.Item(22,3) '192,956,295,255,817'
Set-TextValue $ws.Cells.Item(22,4) '193 Trillion'
Set-TextValue $ws.Cells.Item(22,5) '0.1930% '

# Row 23
$ws.Cells.Item(23,1).Value = 22
Set-TextValue $ws.Cells.Item(23,2) '0x73149b3cd5e1b8536747048259419147e81a71a9'
Set-TextValue $ws.Cells.Item(23,3) '186,576,038,288,042'
Set-TextValue $ws.Cells.Item(23,4) '187 Trillion'
Set-TextValue $ws.Cells.Item(23,5) '0.1866% '

# Row 24
$ws.Cells.Item(24,1).Value = 23
Set-TextValue $ws.Cells.Item(24,2) '0x672c36fa22029369490bb5e33e6d16a7e1309c1e'
Set-TextValue $ws.Cells.Item(24,3) '181,413,020,794,868'
Set-TextValue $ws.Cells.Item(24,4) '181 Trillion'
Set-TextValue $ws.Cells.Item(24,5) '0.1814% '
Set-TextValue $ws.Cells.Item(24,6) '40 Trillion More'
Set-TextValue $ws.Cells.Item(24,7) '40,328,316,152,224'

# Row 25
$ws.Cells.Item(25,1).Value = 24
Set-TextValue $ws.Cells.Item(25,2) '0x497e289791fc2c2b355c259d9516f079d9b52a63'
Set-TextValue $ws.Cells.Item(25,3) '178,884,424,353,453'
Set-TextValue $ws.Cells.Item(25,4) '179 Trillion'
Set-TextValue $ws.Cells.Item(25,5) '0.1789% '

# Row 26
$ws.Cells.Item(26,1).Value = 25
Set-TextValue $ws.Cells.Item(26,2) '0xd693658b31ef059354118ac2b8067b989ebf4b2b'
Set-TextValue $ws.Cells.Item(26,3) '172,179,198,856,685'
Set-TextValue $ws.Cells.Item(26,4) '172 Trillion'
Set-TextValue $ws.Cells.Item(26,5) '0.1722% '
Set-TextValue $ws.Cells.Item(26,6) '5 Trillion More'
Set-TextValue $ws.Cells.Item(26,7) '4,843,670,594,788'

# Row 27
$ws.Cells.Item(27,1).Value = 26
Set-TextValue $ws.Cells.Item(27,2) '0xf3f83f6a5830e55b45b3c44010be0481baa1b9be'
Set-TextValue $ws.Cells.Item(27,3) '170,559,005,253,847'
Set-TextValue $ws.Cells.Item(27,4) '171 Trillion'
Set-TextValue $ws.Cells.Item(27,5) '0.1706% '

# Row 28
$ws.Cells.Item(28,1).Value = 27
Set-TextValue $ws.Cells.Item(28,2) '0x4159fcaefd2216a1b581587ca97da9f53e8ba163'
Set-TextValue $ws.Cells.Item(28,3) '170,029,080,849,270'
Set-TextValue $ws.Cells.Item(28,4) '170 Trillion'
Set-TextValue $ws.Cells.Item(28,5) '0.1700% '

# Row 29
$ws.Cells.Item(29,1).Value = 28
Set-TextValue $ws.Cells.Item(29,2) '0xbdf119001cf9d44d902bf7d8e283e10ab66ddeea'
Set-TextValue $ws.Cells.Item(29,3) '160,512,332,157,587'
Set-TextValue $ws.Cells.Item(29,4) '161 Trillion'
Set-TextValue $ws.Cells.Item(29,5) '0.1605% '

# Row 30
$ws.Cells.Item(30,1).Value = 29
Set-TextValue $ws.Cells.Item(30,2) '0xcc6833974ce5970eac45e7751573c30c7b41a4a5'
Set-TextValue $ws.Cells.Item(30,3) '156,638,608,730,399'
Set-TextValue $ws.Cells.Item(30,4) '157 Trillion'
Set-TextValue $ws.Cells.Item(30,5) '0.1566% '

# Row 31
$ws.Cells.Item(31,1).Value = 30
Set-TextValue $ws.Cells.Item(31,2) '0x7b5b9b8d134bec76023cd6c20358d38714cc5c58'
Set-TextValue $ws.Cells.Item(31,3) '152,751,853,489,352'
Set-TextValue $ws.Cells.Item(31,4) '153 Trillion'
Set-TextValue $ws.Cells.Item(31,5) '0.1528% '

# Row 32
$ws.Cells.Item(32,1).Value = 31
Set-TextValue $ws.Cells.Item(32,2) '0x1ae48253b364374d3db52de311302fc501b87895'
Set-TextValue $ws.Cells.Item(32,3) '152,401,926,728,651'
Set-TextValue $ws.Cells.Item(32,4) '152 Trillion'
Set-TextValue $ws.Cells.Item(32,5) '0.1524% '

# Row 33
$ws.Cells.Item(33,1).Value = 32
Set-TextValue $ws.Cells.Item(33,2) '0x9a7e16cc5d152e60ea52d46d8e422d724bdb4dcf'
Set-TextValue $ws.Cells.Item(33,3) '150,058,675,230,722'
Set-TextValue $ws.Cells.Item(33,4) '150 Trillion'
Set-TextValue $ws.Cells.Item(33,5) '0.1501% '

# Row 44
$ws.Cells.Item(44,1).Value = 43
Set-TextValue $ws.Cells.Item(44,2) '0xde58455ce16cb194a4dc90532326fbf9f3ba8513'
Set-TextValue $ws.Cells.Item(44,3) '102,064,655,233,239'
Set-TextValue $ws.Cells.Item(44,4) '102 Trillion'
Set-TextValue $ws.Cells.Item(44,5) '0.1021% '
Set-TextValue $ws.Cells.Item(44,6) '77 Trillion Less'
Set-TextValue $ws.Cells.Item(44,7) '76,988,621,612,485'

# Row 49
$ws.Cells.Item(49,1).Value = 48
Set-TextValue $ws.Cells.Item(49,2) '0x86437c0875fa78dc98c57bc010ef4ad07bc01715'
Set-TextValue $ws.Cells.Item(49,3) '100,127,887,461,980'
Set-TextValue $ws.Cells.Item(49,4) '100 Trillion'
Set-TextValue $ws.Cells.Item(49,5) '0.1001% '

# Row 53
Set-TextValue $ws.Cells.Item(53,2) 'Total''s'
Set-TextValue $ws.Cells.Item(53,3) '12,045,349,361,697,981'
Set-TextValue $ws.Cells.Item(53,4) '12045 Trillion'
Set-TextValue $ws.Cells.Item(53,5) '12.05%'
Set-TextValue $ws.Cells.Item(53,6) '9 Trillion'
Set-TextValue $ws.Cells.Item(53,7) '8,657,931,074,269'
